# Weekly refresh of the "Hortaliza, Vega Modelo de Temuco - Camote" price
# sheet (commit: "Fruta / hortaliza, semanal").
#
# The published dataset is re-pulled each week; rows keep their market /
# product identity (columns A-C, E-I, R) but the per-offer details -
# Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K-M), Unidad de
# comercializacion (N), Origen (O), Precio $/Kg (P) and Kg o Unidades (Q) -
# get refreshed against the latest weekly pull, which reorders/updates
# those values row by row.
#
# $data holds, per target row number, the new (D,J,K,L,M,N,O,P,Q) tuple.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44498, 20, 20000, 20000, 20000, "`$/malla 20 kilos", "Región de Arica y Parinacota", 1000, 20),
    @(3, 44161, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(4, 44452, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(5, 44466, 20, 25000, 25000, 25000, "`$/caja 15 kilos granel", "Perú", 1667, 15),
    @(6, 44425, 10, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(7, 44448, 45, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(8, 44294, 5, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(9, 44175, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(10, 44364, 15, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(11, 44329, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(12, 44455, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(13, 44369, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(14, 44369, 20, 20000, 20000, 20000, "`$/malla 20 kilos", "Región de Arica y Parinacota", 1000, 20),
    @(15, 44385, 18, 20000, 20000, 20000, "`$/malla 20 kilos", "Región de Arica y Parinacota", 1000, 20),
    @(16, 44341, 40, 17000, 18000, 17500, "`$/malla 20 kilos", "Perú", 875, 20),
    @(17, 44389, 45, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(18, 44321, 15, 25000, 25000, 25000, "`$/caja 15 kilos granel", "Perú", 1667, 15),
    @(19, 44511, 50, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(20, 44186, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(21, 44179, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(22, 44438, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(23, 44188, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(24, 44441, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(26, 44496, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(27, 44315, 30, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(28, 44315, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Región de Arica y Parinacota", 1000, 20),
    @(29, 44340, 40, 18000, 18000, 18000, "`$/malla 20 kilos", "Perú", 900, 20),
    @(30, 44497, 30, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(31, 44497, 40, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(32, 44525, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(33, 44508, 40, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Perú", 1333, 15),
    @(34, 44316, 20, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15),
    @(35, 44512, 30, 20000, 20000, 20000, "`$/malla 20 kilos", "Perú", 1000, 20),
    @(36, 44424, 30, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Región de Arica y Parinacota", 1333, 15)
)

foreach ($row in $data) {
    $r  = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[6]   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $row[7]   # O: Origen
    $ws.Cells.Item($r, 16).Value = $row[8]   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $row[9]   # Q: Kg o Unidades
}
